# Add a new "Italy" market sheet, based on the existing "Slovakia" sheet,
# matching the layout/content of the other per-market accessory sheets.

$wb = $excel.ActiveWorkbook

# The Slovakia sheet is the right-most tab and serves as the template for
# the new Italy sheet (same columns/rows/styles, different market name and
# reference code).
$src = $wb.Worksheets.Item("Slovakia")

# Duplicate it (places the copy immediately after the source sheet) - this
# mirrors "Move or Copy... (Create a copy)" in the Excel UI.
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Italy"

# Update the market name / reference code cells for Italy.
$newSheet.Range("B2").Value = "Italy Market"

# B4 loses its inherited fill/border formatting (typed fresh rather than
# pasted), so clear formatting before writing the new value.
$newSheet.Range("B4").ClearFormats()
$newSheet.Range("B4").Value = "NGC-3145/T2153/T2218/T2456"

# Append two extra rows (18 & 19) that repeat the final "Wg" / "Accessories"
# rows (16 & 17), keeping their original styling.
$newSheet.Range("A16:A17").Copy()
$newSheet.Range("A18:A19").PasteSpecial(-4122)  # xlPasteFormats
$wgValue = $newSheet.Range("A16").Value2
$accessoriesValue = $newSheet.Range("A17").Value2
$newSheet.Range("A18").Value2 = $wgValue
$newSheet.Range("A19").Value2 = $accessoriesValue

# The source (Slovakia) sheet ends up with everything selected and is no
# longer the active tab.
$src.Cells.Select()

# The new Italy sheet becomes the active tab, with B18 selected.
$newSheet.Activate()
$newSheet.Range("B18").Select()
